$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 95 ---
$ws.Cells.Item(95, 1).Value = 45474.2916666667
$ws.Cells.Item(95, 2).Value = 0
$ws.Cells.Item(95, 3).Value = 6.17999982833862
$ws.Cells.Item(95, 4).Value = 6.17999982833862
$ws.Cells.Item(95, 5).Value = 6.17999982833862
$ws.Cells.Item(95, 6).Value = 6.17999982833862

# G95 must be stored as text (shared string) even though it looks numeric
$ws.Cells.Item(95, 7).NumberFormat = "@"
$ws.Cells.Item(95, 7).Value = "6.17999982833862"
$ws.Cells.Item(95, 7).Style = "Normal"

$ws.Cells.Item(95, 8).Value = "PAL.MI"

# --- Row 96 ---
$ws.Cells.Item(96, 1).Value = 45475.6447685185
$ws.Cells.Item(96, 2).Value = 2400
$ws.Cells.Item(96, 3).Value = 6.26000022888184
$ws.Cells.Item(96, 4).Value = 6.17999982833862
$ws.Cells.Item(96, 5).Value = 6.23999977111816
$ws.Cells.Item(96, 6).Value = 6.23999977111816

# G96 must be stored as text (shared string) even though it looks numeric
$ws.Cells.Item(96, 7).NumberFormat = "@"
$ws.Cells.Item(96, 7).Value = "6.23999977111816"
$ws.Cells.Item(96, 7).Style = "Normal"

$ws.Cells.Item(96, 8).Value = "PAL.MI"

# Match the date style (style index 1, yyyy-mm-dd hh:mm:ss) used by column A elsewhere
$ws.Range("A94").Copy()
$ws.Range("A95:A96").PasteSpecial(-4122)
